{"js": "// Apply the LOB1233 syllabus updates:\n// 1) Bump the \"Ativa\u00e7\u00e3o\" date from 2020 to 2025.\n// 2) Append \"; Gest\u00e3o de Recursos H\u00eddricos.\" to the PT \"Programa resumido\" text.\n// 3) Append \", Water Resources Management.\" to the EN \"Programa resumido\" text.\n// 4) Append extra sentences to the PT \"Programa\" text.\n// 5) Append extra sentences to the EN \"Programa\" text.\n\nconst replacements = [\n  {\n    find: \"Ativa\u00e7\u00e3o: 01/01/2020\",\n    replace: \"Ativa\u00e7\u00e3o: 01/01/2025\",\n  },\n  {\n    find:\n      \"Formas de representa\u00e7\u00e3o e apresenta\u00e7\u00e3o das caracter\u00edsticas ambientais; Caracteriza\u00e7\u00e3o ambiental e sua aplica\u00e7\u00e3o em uma bacia hidrogr\u00e1fica; Determina\u00e7\u00e3o das suscetibilidades e voca\u00e7\u00f5es do meio ambiente e o conceito de sustentabilidade ambiental.\",\n    replace:\n      \"Formas de representa\u00e7\u00e3o e apresenta\u00e7\u00e3o das caracter\u00edsticas ambientais; Caracteriza\u00e7\u00e3o ambiental e sua aplica\u00e7\u00e3o em uma bacia hidrogr\u00e1fica; Determina\u00e7\u00e3o das suscetibilidades e voca\u00e7\u00f5es do meio ambiente e o conceito de sustentabilidade ambiental; Gest\u00e3o de Recursos H\u00eddricos.\",\n  },\n  {\n    find:\n      \"Representation and presentation forms of environmental characteristics; Environmental characterization and its application in a watershed; Environment susceptibilities and vocations determination and environmental susceptibility concept.\",\n    replace:\n      \"Representation and presentation forms of environmental characteristics; Environmental characterization and its application in a watershed; Environment susceptibilities and vocations determination and environmental susceptibility concept, Water Resources Management.\",\n  },\n  {\n    find:\n      \"Bacia hidrogr\u00e1fica como unidade de estudo e gest\u00e3o de recursos h\u00eddricos; caracteriza\u00e7\u00e3o ambiental de bacias hidrogr\u00e1ficas; Caracteriza\u00e7\u00e3o morfom\u00e9trica de bacias hidrogr\u00e1ficas; Regi\u00f5es Hidrogr\u00e1ficas do Brasil; Estrutura institucional e marcos legais em recursos h\u00eddricos no Brasil.\",\n    replace:\n      \"Bacia hidrogr\u00e1fica como unidade de estudo e gest\u00e3o de recursos h\u00eddricos; caracteriza\u00e7\u00e3o ambiental de bacias hidrogr\u00e1ficas; Caracteriza\u00e7\u00e3o morfom\u00e9trica de bacias hidrogr\u00e1ficas; Regi\u00f5es Hidrogr\u00e1ficas do Brasil; Estrutura institucional e marcos legais em recursos h\u00eddricos no Brasil. Legisla\u00e7\u00e3o e instrumentos pertinentes. A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.\",\n  },\n  {\n    find:\n      \"Watershed as a unit of study and management of water resources; Environmental characterization of watershed; Morphometric characterization of watershed; Brazilian hydrographic Regions; Brazilian institutional structure and legal frameworks in water resources.\",\n    replace:\n      \"Watershed as a unit of study and management of water resources; Environmental characterization of watershed; Morphometric characterization of watershed; Brazilian hydrographic Regions; Brazilian institutional structure and legal frameworks in water resources; Legislation and relevant instrument. The discipline may have didactic trips to complement the content of the discipline.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the LOB1233 syllabus updates:\n# 1) Bump the \"Ativa\u00e7\u00e3o\" date from 2020 to 2025.\n# 2) Append \"; Gest\u00e3o de Recursos H\u00eddricos.\" to the PT \"Programa resumido\" text.\n# 3) Append \", Water Resources Management.\" to the EN \"Programa resumido\" text.\n# 4) Append extra sentences to the PT \"Programa\" text.\n# 5) Append extra sentences to the EN \"Programa\" text.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute(\n        $findText,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $replaceText,# ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2020\" \"Ativa\u00e7\u00e3o: 01/01/2025\"\n\nReplace-Text \"Formas de representa\u00e7\u00e3o e apresenta\u00e7\u00e3o das caracter\u00edsticas ambientais; Caracteriza\u00e7\u00e3o ambiental e sua aplica\u00e7\u00e3o em uma bacia hidrogr\u00e1fica; Determina\u00e7\u00e3o das suscetibilidades e voca\u00e7\u00f5es do meio ambiente e o conceito de sustentabilidade ambiental.\" \"Formas de representa\u00e7\u00e3o e apresenta\u00e7\u00e3o das caracter\u00edsticas ambientais; Caracteriza\u00e7\u00e3o ambiental e sua aplica\u00e7\u00e3o em uma bacia hidrogr\u00e1fica; Determina\u00e7\u00e3o das suscetibilidades e voca\u00e7\u00f5es do meio ambiente e o conceito de sustentabilidade ambiental; Gest\u00e3o de Recursos H\u00eddricos.\"\n\nReplace-Text \"Representation and presentation forms of environmental characteristics; Environmental characterization and its application in a watershed; Environment susceptibilities and vocations determination and environmental susceptibility concept.\" \"Representation and presentation forms of environmental characteristics; Environmental characterization and its application in a watershed; Environment susceptibilities and vocations determination and environmental susceptibility concept, Water Resources Management.\"\n\nReplace-Text \"Bacia hidrogr\u00e1fica como unidade de estudo e gest\u00e3o de recursos h\u00eddricos; caracteriza\u00e7\u00e3o ambiental de bacias hidrogr\u00e1ficas; Caracteriza\u00e7\u00e3o morfom\u00e9trica de bacias hidrogr\u00e1ficas; Regi\u00f5es Hidrogr\u00e1ficas do Brasil; Estrutura institucional e marcos legais em recursos h\u00eddricos no Brasil.\" \"Bacia hidrogr\u00e1fica como unidade de estudo e gest\u00e3o de recursos h\u00eddricos; caracteriza\u00e7\u00e3o ambiental de bacias hidrogr\u00e1ficas; Caracteriza\u00e7\u00e3o morfom\u00e9trica de bacias hidrogr\u00e1ficas; Regi\u00f5es Hidrogr\u00e1ficas do Brasil; Estrutura institucional e marcos legais em recursos h\u00eddricos no Brasil. Legisla\u00e7\u00e3o e instrumentos pertinentes. A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.\"\n\nReplace-Text \"Watershed as a unit of study and management of water resources; Environmental characterization of watershed; Morphometric characterization of watershed; Brazilian hydrographic Regions; Brazilian institutional structure and legal frameworks in water resources.\" \"Watershed as a unit of study and management of water resources; Environmental characterization of watershed; Morphometric characterization of watershed; Brazilian hydrographic Regions; Brazilian institutional structure and legal frameworks in water resources; Legislation and relevant instrument. The discipline may have didactic trips to complement the content of the discipline.\"\n"}
